# Update the "MODALIDAD" related district-level info (AT_ETL_102022 / AT_GORE_102022 / AT_ETR_102022 columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - AYACUCHO
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 1

# Row 9 - CUSCO
$ws.Range("F9").Value = 1

# Row 12 - ICA
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 1

# Row 15 - LAMBAYEQUE
$ws.Range("G15").Value = 2

# Row 24 - SAN MARTIN
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 2

# Row 25 - TACNA
$ws.Range("F25").Value = 2

# Row 27 - UCAYALI
$ws.Range("F27").Value = 2
